# Finally pos tagging the whole dataset.
# Adds a new column G ("Brown+Wordnet+Names+COCA") of POS tags to Sheet1,
# mirroring the existing B-F "tagger" columns, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header (row 1) -------------------------------------------------
$headerCell = $ws.Range("G1")
$headerCell.Value = "Brown+Wordnet+Names+COCA"
$headerCell.Font.Bold = $true

# --- Data rows (2-79) ------------------------------------------------
$rows = @(
    @{ Row = 2; Value = "NN"; Style = "Bad" },
    @{ Row = 3; Value = "NN"; Style = "Bad" },
    @{ Row = 4; Value = "NNP"; Style = "Good" },
    @{ Row = 5; Value = "NN"; Style = "Good" },
    @{ Row = 6; Value = "NN"; Style = "Good" },
    @{ Row = 7; Value = "NN"; Style = "Good" },
    @{ Row = 8; Value = "NN"; Style = "Good" },
    @{ Row = 9; Value = "NN"; Style = "Good" },
    @{ Row = 10; Value = "NN"; Style = "Good" },
    @{ Row = 11; Value = "UH"; Style = "Good" },
    @{ Row = 12; Value = "RB"; Style = "Good" },
    @{ Row = 13; Value = "JJ"; Style = "Good" },
    @{ Row = 14; Value = "NN"; Style = "Good" },
    @{ Row = 15; Value = "NN"; Style = "Good" },
    @{ Row = 16; Value = "NNP"; Style = "Good" },
    @{ Row = 17; Value = "JJ"; Style = "Good" },
    @{ Row = 18; Value = "NN"; Style = "Good" },
    @{ Row = 19; Value = "NN"; Style = "Good" },
    @{ Row = 20; Value = "NN"; Style = "Good" },
    @{ Row = 21; Value = "VB"; Style = "Good" },
    @{ Row = 22; Value = "RP"; Style = "Good" },
    @{ Row = 23; Value = "NN"; Style = "Good" },
    @{ Row = 24; Value = "NN"; Style = "Good" },
    @{ Row = 25; Value = "NN"; Style = "Bad" },
    @{ Row = 26; Value = "PPO"; Style = "Good" },
    @{ Row = 27; Value = "JJ"; Style = "Good" },
    @{ Row = 28; Value = "NN"; Style = "Good" },
    @{ Row = 29; Value = "NN"; Style = "Bad" },
    @{ Row = 30; Value = "KK"; Style = "Neutral" },
    @{ Row = 31; Value = "NN"; Style = "Good" },
    @{ Row = 32; Value = "JJ"; Style = "Good" },
    @{ Row = 33; Value = "NN"; Style = "Good" },
    @{ Row = 34; Value = "NN"; Style = "Bad" },
    @{ Row = 35; Value = "NN"; Style = "Bad" },
    @{ Row = 36; Value = "PPSS"; Style = "Good" },
    @{ Row = 37; Value = "NN"; Style = "Good" },
    @{ Row = 38; Value = "IN"; Style = "Good" },
    @{ Row = 39; Value = "NN"; Style = "Good" },
    @{ Row = 40; Value = "NN"; Style = "Good" },
    @{ Row = 41; Value = "VB"; Style = "Good" },
    @{ Row = 42; Value = "NN"; Style = "Bad" },
    @{ Row = 43; Value = "VB"; Style = "Good" },
    @{ Row = 44; Value = "PPO"; Style = "Good" },
    @{ Row = 45; Value = "NNP"; Style = "Good" },
    @{ Row = 46; Value = "NP"; Style = "Good" },
    @{ Row = 47; Value = "NNP"; Style = "Good" },
    @{ Row = 48; Value = "NP"; Style = "Good" },
    @{ Row = 49; Value = "NN"; Style = "Good" },
    @{ Row = 50; Value = "NN"; Style = "Good" },
    @{ Row = 51; Value = "NP"; Style = "Good" },
    @{ Row = 52; Value = "NP"; Style = "Good" },
    @{ Row = 53; Value = "NNP"; Style = "Good" },
    @{ Row = 54; Value = "NN"; Style = "Bad" },
    @{ Row = 55; Value = "NP"; Style = "Good" },
    @{ Row = 56; Value = "VB"; Style = "Good" },
    @{ Row = 57; Value = "NNS"; Style = "Good" },
    @{ Row = 58; Value = "NN"; Style = "Good" },
    @{ Row = 59; Value = "JJ"; Style = "Good" },
    @{ Row = 60; Value = "NN"; Style = "Bad" },
    @{ Row = 61; Value = "KK"; Style = "Neutral" },
    @{ Row = 62; Value = "VB"; Style = "Good" },
    @{ Row = 63; Value = "PPL"; Style = "Good" },
    @{ Row = 64; Value = "NN"; Style = "Good" },
    @{ Row = 65; Value = "NN"; Style = "Good" },
    @{ Row = 66; Value = "NN"; Style = "Good" },
    @{ Row = 67; Value = "NN"; Style = "Good" },
    @{ Row = 68; Value = "NN"; Style = "Good" },
    @{ Row = 69; Value = "NN"; Style = "Good" },
    @{ Row = 70; Value = "NN"; Style = "Good" },
    @{ Row = 71; Value = "VB"; Style = "Good" },
    @{ Row = 72; Value = "AP"; Style = "Good" },
    @{ Row = 73; Value = "NNS"; Style = "Good" },
    @{ Row = 74; Value = "AT"; Style = "Good" },
    @{ Row = 75; Value = "NN"; Style = "Good" },
    @{ Row = 76; Value = "IN"; Style = "Good" },
    @{ Row = 77; Value = "NN"; Style = "Good" },
    @{ Row = 78; Value = "CC"; Style = "Good" },
    @{ Row = 79; Value = "NN"; Style = "Good" }
)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r.Row, 7)
    $cell.Value = $r.Value
    $cell.Style = $r.Style
}

# --- Selection reflects the newly added column ------------------------
$ws.Range("G1").Select()
